$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: No Route Table data (red font) - was "No Route Table Data"
$ws.Range("C3").Value2 = "No Route Table data"
$ws.Range("C3").Font.Color = 255

# Row 5: VPC - was "Hosting VPC"
$ws.Range("C5").Value2 = "VPC"

# Row 7: Default Route Table' (quote-prefixed) - was curly-quoted ‘Default Route Table’
$ws.Range("C7").Value2 = "'Default Route Table'"

# Row 8: Customized Route Table' (quote-prefixed) - was curly-quoted ‘Customized Table’
$ws.Range("C8").Value2 = "'Customized Route Table'"

# Row 16: Pop Layer Operations - text unchanged, now flagged red for QA
$ws.Range("C16").Font.Color = 255

# Row 17: Begin to get Route Table list (red) - was curly-quoted ‘Begin to get Route Table List’
$ws.Range("C17").Value2 = "Begin to get Route Table list"
$ws.Range("C17").Font.Color = 255

# Row 18: Route Table list data are (red) - was curly-quoted ‘Route Table List data are
$ws.Range("C18").Value2 = "Route Table list data are"
$ws.Range("C18").Font.Color = 255

# Selection moves to B16
$ws.Range("B16").Select()

# Explicit (default) page setup, portrait orientation
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 0
